$wb = $excel.ActiveWorkbook

# ----- Sheet: ALC -----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H70").Value = 1893.8889
$ws.Range("I70").Value = 1233.3334
$ws.Range("J70").Value = 2224.1667
$ws.Range("K70").Value = 3700.0002
$ws.Range("L70").Value = 6672.500100000001
$ws.Range("M70").Value = -3430.0002
$ws.Range("N70").Value = -7212.500100000001
$ws.Range("H73").Value = 1893.8889
$ws.Range("I73").Value = 1233.3334
$ws.Range("J73").Value = 2224.1667
$ws.Range("K73").Value = 3700.0002
$ws.Range("L73").Value = 6672.500100000001
$ws.Range("M73").Value = -2764.0002
$ws.Range("N73").Value = -8544.500100000001
$ws.Range("H105").Value = 34835.5
$ws.Range("J105").Value = 34835.5
$ws.Range("L105").Value = 34835.5
$ws.Range("N105").Value = -41823.5

# ----- Sheet: ARM -----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 256.55554
$ws.Range("I4").Value = 46.5
$ws.Range("J4").Value = 676.6667
$ws.Range("K4").Value = 46.5
$ws.Range("L4").Value = 676.6667
$ws.Range("M4").Value = 69.5
$ws.Range("N4").Value = -908.6667
$ws.Range("H27").Value = 3032.2666
$ws.Range("I27").Value = 0
$ws.Range("J27").Value = 3032.2666
$ws.Range("K27").Value = 0
$ws.Range("L27").Value = 3032.2666
$ws.Range("M27").ClearContents()
$ws.Range("N27").Value = -3400.2666
$ws.Range("H32").Value = 500699.25
$ws.Range("I32").Value = 4788.6
$ws.Range("J32").Value = 2217313
$ws.Range("K32").Value = 4788.6
$ws.Range("L32").Value = 2217313
$ws.Range("M32").Value = -4501.6
$ws.Range("N32").Value = -2217887
$ws.Range("H43").Value = 11062.667
$ws.Range("I43").Value = 4999
$ws.Range("J43").Value = 12275.4
$ws.Range("K43").Value = 4999
$ws.Range("L43").Value = 12275.4
$ws.Range("M43").Value = -4686
$ws.Range("N43").Value = -12901.4
$ws.Range("H45").Value = 4540.24
$ws.Range("I45").Value = 4639.385
$ws.Range("J45").Value = 4432.8335
$ws.Range("K45").Value = 4639.385
$ws.Range("L45").Value = 4432.8335
$ws.Range("M45").Value = -4262.385
$ws.Range("N45").Value = -5186.8335
$ws.Range("H61").Value = 2354.08
$ws.Range("I61").Value = 2269.2727
$ws.Range("J61").Value = 2976
$ws.Range("K61").Value = 2269.2727
$ws.Range("L61").Value = 2976
$ws.Range("M61").Value = -2057.2727
$ws.Range("N61").Value = -3400
$ws.Range("H101").Value = 0
$ws.Range("J101").Value = 0
$ws.Range("L101").Value = 0
$ws.Range("N101").ClearContents()
$ws.Range("H104").Value = 76408.336
$ws.Range("J104").Value = 76408.336
$ws.Range("L104").Value = 76408.336
$ws.Range("N104").Value = -83396.336
$ws.Range("H106").Value = 0
$ws.Range("J106").Value = 0
$ws.Range("L106").Value = 0
$ws.Range("N106").ClearContents()
$ws.Range("H136").Value = 2354.08
$ws.Range("I136").Value = 2269.2727
$ws.Range("J136").Value = 2976
$ws.Range("K136").Value = 6807.8181
$ws.Range("L136").Value = 8928
$ws.Range("M136").Value = -4257.8181
$ws.Range("N136").Value = -14028

# ----- Sheet: BSM -----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H92").Value = 20200.5
$ws.Range("J92").Value = 20200.5
$ws.Range("L92").Value = 20200.5
$ws.Range("N92").Value = -25192.5
$ws.Range("H100").Value = 24900
$ws.Range("J100").Value = 24900
$ws.Range("L100").Value = 24900
$ws.Range("N100").Value = -27064

# ----- Sheet: CRP -----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1988.8823
$ws.Range("I16").Value = 1557.9286
$ws.Range("J16").Value = 4000
$ws.Range("K16").Value = 1557.9286
$ws.Range("L16").Value = 4000
$ws.Range("M16").Value = -1270.9286
$ws.Range("N16").Value = -4574
$ws.Range("H88").Value = 13000
$ws.Range("J88").Value = 18000
$ws.Range("L88").Value = 18000
$ws.Range("N88").Value = -18812
$ws.Range("H91").Value = 13000
$ws.Range("J91").Value = 18000
$ws.Range("L91").Value = 18000
$ws.Range("N91").Value = -20808
$ws.Range("H113").Value = 1988.8823
$ws.Range("I113").Value = 1557.9286
$ws.Range("J113").Value = 4000
$ws.Range("K113").Value = 1557.9286
$ws.Range("L113").Value = 4000
$ws.Range("M113").Value = 612.0714
$ws.Range("N113").Value = -8340
$ws.Range("H134").Value = 150018500
$ws.Range("I134").Value = 200001040
$ws.Range("J134").Value = 70900
$ws.Range("K134").Value = 600003120
$ws.Range("L134").Value = 212700
$ws.Range("M134").Value = -600000585
$ws.Range("N134").Value = -217770

# ----- Sheet: CUL -----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H9").Value = 971.4286
$ws.Range("I9").Value = 200
$ws.Range("J9").Value = 2000
$ws.Range("K9").Value = 600
$ws.Range("L9").Value = 6000
$ws.Range("M9").Value = -376
$ws.Range("N9").Value = -6448
$ws.Range("H10").Value = 0
$ws.Range("I10").Value = 0
$ws.Range("K10").Value = 0
$ws.Range("M10").ClearContents()

# ----- Sheet: GSM -----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 35.727272
$ws.Range("I2").Value = 16.470589
$ws.Range("J2").Value = 101.2
$ws.Range("K2").Value = 16.470589
$ws.Range("L2").Value = 101.2
$ws.Range("M2").Value = 96.529411
$ws.Range("N2").Value = -327.2
$ws.Range("H12").Value = 90003
$ws.Range("I12").Value = 90003
$ws.Range("K12").Value = 90003
$ws.Range("M12").Value = -89863
$ws.Range("H14").Value = 30002380
$ws.Range("I14").Value = 37500476
$ws.Range("K14").Value = 37500476
$ws.Range("M14").Value = -37500308
$ws.Range("H15").Value = 15000
$ws.Range("J15").Value = 15000
$ws.Range("L15").Value = 15000
$ws.Range("N15").Value = -15576
$ws.Range("H44").Value = 400978
$ws.Range("J44").Value = 400978
$ws.Range("L44").Value = 400978
$ws.Range("N44").Value = -402170
$ws.Range("H81").Value = 15000
$ws.Range("J81").Value = 15000
$ws.Range("L81").Value = 15000
$ws.Range("N81").Value = -16996
$ws.Range("H84").Value = 15000
$ws.Range("J84").Value = 15000
$ws.Range("L84").Value = 45000
$ws.Range("N84").Value = -54984
$ws.Range("H95").Value = 15000
$ws.Range("J95").Value = 15000
$ws.Range("L95").Value = 15000
$ws.Range("N95").Value = -20492
$ws.Range("H98").Value = 14680.857
$ws.Range("J98").Value = 14680.857
$ws.Range("L98").Value = 14680.857
$ws.Range("N98").Value = -20670.857
$ws.Range("H101").Value = 37200
$ws.Range("J101").Value = 37200
$ws.Range("L101").Value = 37200
$ws.Range("N101").Value = -43690

# ----- Sheet: LTW -----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H51").Value = 0
$ws.Range("J51").Value = 0
$ws.Range("L51").Value = 0
$ws.Range("N51").ClearContents()

# ----- Sheet: WVR -----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H41").Value = 21733
$ws.Range("J41").Value = 12021.833
$ws.Range("L41").Value = 12021.833
$ws.Range("N41").Value = -12801.833
$ws.Range("H45").Value = 9335.929
$ws.Range("J45").Value = 9335.929
$ws.Range("L45").Value = 9335.929
$ws.Range("N45").Value = -10317.929
$ws.Range("H74").Value = 12193.5
$ws.Range("J74").Value = 12193.5
$ws.Range("L74").Value = 12193.5
$ws.Range("N74").Value = -14065.5
$ws.Range("H77").Value = 12193.5
$ws.Range("J77").Value = 12193.5
$ws.Range("L77").Value = 36580.5
$ws.Range("N77").Value = -45940.5
$ws.Range("H126").Value = 2020.1428
$ws.Range("I126").Value = 2079.7273
$ws.Range("J126").Value = 1801.6666
$ws.Range("K126").Value = 6239.1819
$ws.Range("L126").Value = 5404.9998
$ws.Range("M126").Value = -3769.1819
$ws.Range("N126").Value = -10344.9998
$ws.Range("H136").Value = 585.6896400000001
$ws.Range("I136").Value = 319.2
$ws.Range("J136").Value = 2251.25
$ws.Range("K136").Value = 957.5999999999999
$ws.Range("L136").Value = 6753.75
$ws.Range("M136").Value = 1592.4
$ws.Range("N136").Value = -11853.75
